$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

$wsSchedule.Range("E2").Value = 518.792508
$wsSchedule.Range("F2").Value = 34.31167380952381
$wsSchedule.Range("E3").Value = -233.10969825
$wsSchedule.Range("F3").Value = -7.70865404265873
$wsSchedule.Range("E4").Value = 493.37327325
$wsSchedule.Range("F4").Value = 32.63050749007937
$wsDetailed.Range("B7").Value = 64.89
$wsDetailed.Range("B8").Value = 65
$wsDetailed.Range("B9").Value = 78
$wsDetailed.Range("B10").Value = 84.58459000000001
$wsDetailed.Range("C10").Value = "historical"
$wsDetailed.Range("B11").Value = 78
$wsDetailed.Range("B12").Value = 83.48896999999999
$wsDetailed.Range("B13").Value = 105.79
$wsDetailed.Range("B14").Value = 105
$wsDetailed.Range("B15").Value = 82.98098
$wsDetailed.Range("B16").Value = 56.98
$wsDetailed.Range("B17").Value = 9.2277
$wsDetailed.Range("B18").Value = 0.009560000000000001
$wsDetailed.Range("B19").Value = -5.74313
$wsDetailed.Range("B20").Value = -6.57611
$wsDetailed.Range("B21").Value = -6.91994
$wsDetailed.Range("B22").Value = -8.29374
$wsDetailed.Range("B23").Value = -9.27102
$wsDetailed.Range("B24").Value = -14
$wsDetailed.Range("B25").Value = -15.1557
$wsDetailed.Range("B26").Value = -15.42766
$wsDetailed.Range("B27").Value = -20
$wsDetailed.Range("B28").Value = -19.98
$wsDetailed.Range("B29").Value = -22.3004
$wsDetailed.Range("B30").Value = -24.42766
$wsDetailed.Range("B31").Value = -23.73604
$wsDetailed.Range("B32").Value = -24.12776
$wsDetailed.Range("B33").Value = -23.13727
$wsDetailed.Range("B34").Value = -6.88088
$wsDetailed.Range("B36").Value = 36.06018
$wsDetailed.Range("B37").Value = 49.47591
$wsDetailed.Range("B38").Value = 56.35715
$wsDetailed.Range("B39").Value = 64.46167
$wsDetailed.Range("B40").Value = 73.27
$wsDetailed.Range("B44").Value = 57.46248
$wsDetailed.Range("B45").Value = 60.00504
$wsDetailed.Range("B46").Value = 57.31
$wsDetailed.Range("B47").Value = 63.17615
$wsDetailed.Range("B48").Value = 64.8901
$wsDetailed.Range("B49").Value = 64.8901
